$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "310.78"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.49%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.37%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07766"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.83%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.52%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.223"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.43%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.876"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-7.64%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9205"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.96%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1210"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-5.33%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1899"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.45%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09221"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "5.36%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03429"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.72%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09681"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.68%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.01%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005942"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.56%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.560"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.84%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.056"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.85%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.06%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.261"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "4.94%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1269"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.93%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.02106"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "5,593.12%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04358"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.41%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001199"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.70%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004252"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-8.10%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02093"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-7.73%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05027"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.01%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007700"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.64%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009805"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.99%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1348"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.40%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002182"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "4.14%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009572"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "8.37%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006715"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "3.37%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.16%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.001201"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.18%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002937"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-2.29%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.16%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.16%"
